$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column P (Unemployment Rate) values that were re-rounded to 1 decimal place
$ws.Range("P14").Value = 8.8
$ws.Range("P16").Value = 7.3
$ws.Range("P17").Value = 6.7
$ws.Range("P20").Value = 5.4
$ws.Range("P24").Value = 4.2
$ws.Range("P26").Value = 4
$ws.Range("P28").Value = 4.4
$ws.Range("P30").Value = 4.5
$ws.Range("P31").Value = 4.6
$ws.Range("P33").Value = 5.3
$ws.Range("P34").Value = 4.8
$ws.Range("P35").Value = 4.8
$ws.Range("P37").Value = 4.8
$ws.Range("P38").Value = 5.1
$ws.Range("P39").Value = 4.7
$ws.Range("P42").Value = 4.4
$ws.Range("P45").Value = 4.5
$ws.Range("P47").Value = 4.8
$ws.Range("P49").Value = 4.6
$ws.Range("P53").Value = 5.2
$ws.Range("P54").Value = 5.4
$ws.Range("P55").Value = 5.9
$ws.Range("P56").Value = 7.4
$ws.Range("P57").Value = 8.6
$ws.Range("P58").Value = 10.8
$ws.Range("P61").Value = 13.8
$ws.Range("P63").Value = 14.4
$ws.Range("P65").Value = 15.6
$ws.Range("P66").Value = 15.2
$ws.Range("P67").Value = 15.1
$ws.Range("P69").Value = 15.9
$ws.Range("P70").Value = 15.9
$ws.Range("P71").Value = 15.5
$ws.Range("P72").Value = 15.5
$ws.Range("P74").Value = 14.5
$ws.Range("P75").Value = 14.4
$ws.Range("P76").Value = 13.3
$ws.Range("P77").Value = 12.9
$ws.Range("P78").Value = 12.8
$ws.Range("P80").Value = 11.7
$ws.Range("P81").Value = 10.9
$ws.Range("P82").Value = 10.6
$ws.Range("P84").Value = 9.6
$ws.Range("P85").Value = 9.4
$ws.Range("P86").Value = 9
$ws.Range("P87").Value = 8.8
$ws.Range("P88").Value = 8.2
$ws.Range("P89").Value = 7.5
$ws.Range("P90").Value = 7.3
$ws.Range("P91").Value = 6.7
$ws.Range("P92").Value = 6.7
$ws.Range("P93").Value = 6.4
$ws.Range("P94").Value = 6
$ws.Range("P95").Value = 5.9
$ws.Range("P97").Value = 5.7
$ws.Range("P99").Value = 5.2
$ws.Range("P101").Value = 4.8
$ws.Range("P102").Value = 4.9
$ws.Range("P103").Value = 5.1
$ws.Range("P105").Value = 6.3
$ws.Range("P106").Value = 7.3
$ws.Range("P107").Value = 7
$ws.Range("P108").Value = 5.5
$ws.Range("P109").Value = 5.2
$ws.Range("P110").Value = 4.9
$ws.Range("P112").Value = 4.2
$ws.Range("P113").Value = 4.4
$ws.Range("P115").Value = 4.2
$ws.Range("P117").Value = 4.5
$ws.Range("P120").Value = 4.2

# Add new Payroll employment value for row 123
$ws.Range("N123").Value = 2818.1
